# Regenerate quadratic/linear problem data (mirrors commit:
# "volver a generar problemas cuadraticos y lineales").
#
# The source tool regenerates the whole workbook's cell contents from
# scratch, so every cell (including the ones whose value doesn't
# actually change) is rewritten here, sheet by sheet, top-left to
# bottom-right -- this reproduces the shared-string table order that a
# from-scratch regeneration produces.

function Set-TextValue {
    # Write $text into $range as TEXT even when it "looks" numeric
    # (otherwise Excel auto-converts e.g. "-9.25" into a real number).
    # Mark the cell as Text, assign, then strip the format back off so
    # no visible number-format change remains on the cell itself.
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# Wipe every sheet's contents first so the regenerated shared-string
# table starts clean (old, no-longer-referenced strings drop out).
for ($i = 1; $i -le 7; $i++) {
    $wb.Worksheets.Item($i).Cells.Clear()
}

# ---------------------------------------------------------------
# Funciones_Objetivo
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "Leader_Expr"
$ws.Range("B1").Value = "Follower_Expr"
$ws.Range("A2").Value = "(-3 + x)^2 + (-2 + y)^2"
$ws.Range("B2").Value = "(-5 + y)^2"

# ---------------------------------------------------------------
# Restricciones_del_lider
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Range("A1").Value = "Expression"
$ws.Range("B1").Value = "Function_Evaluation"
$ws.Range("C1").Value = "Restriction_Set_Type"
$ws.Range("D1").Value = "MIU_value"

# ---------------------------------------------------------------
# Restricciones_del_follower
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Range("A1").Value = "Expression"
$ws.Range("B1").Value = "Function_Evaluation"
$ws.Range("C1").Value = "Restriction_Set_Type"
$ws.Range("D1").Value = "Lambda_value"
$ws.Range("E1").Value = "Beta_value"
$ws.Range("F1").Value = "Gamma_value"

$ws.Range("A2").Value = "-4.75 + x"
Set-TextValue $ws.Range("B2") "-9.25"
$ws.Range("C2").Value = "J_0_L0_v"
Set-TextValue $ws.Range("D2") "0.07"
Set-TextValue $ws.Range("E2") "1.0"
Set-TextValue $ws.Range("F2") "0"

$ws.Range("A3").Value = "-13.902999999999999 + x + 2.26y"
Set-TextValue $ws.Range("B3") "11.902999999999999"
$ws.Range("C3").Value = "J_0_LP_v"
Set-TextValue $ws.Range("D3") "0.21"
Set-TextValue $ws.Range("E3") "-2.3000000000000003"
Set-TextValue $ws.Range("F3") "-5.6000000000000005"

$ws.Range("A4").Value = "-28.802816901408452 - 2x - 4.272300469483569y"
Set-TextValue $ws.Range("B4") "-27.802816901408452"
$ws.Range("C4").Value = "J_Ne_L0_v"
Set-TextValue $ws.Range("D4") "0.4"
Set-TextValue $ws.Range("E4") "-1.0"
Set-TextValue $ws.Range("F4") "-9.1"

# ---------------------------------------------------------------
# Punto_modificado
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Range("A1").Value = "x"
$ws.Range("B1").Value = "y"
Set-TextValue $ws.Range("A2") "4.75"
Set-TextValue $ws.Range("B2") "4.05"

# ---------------------------------------------------------------
# Vector_bf
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
$ws.Range("A1").Value = "vec_bf"
Set-TextValue $ws.Range("A2") "3.134320187793428"

# ---------------------------------------------------------------
# Vector_BF
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
$ws.Range("A1").Value = "vec_BF"
Set-TextValue $ws.Range("A2") "-4.199999999999999"
Set-TextValue $ws.Range("A3") "-7.434300469483568"

# ---------------------------------------------------------------
# Vector_Alpha (A2 here really is a numeric cell, not text)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(7)
$ws.Range("A1").Value = "vec_alpha"
$ws.Range("A2").Value = 2.13
